$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -18.54418339531084
$ws.Range("C2").Value = 2.012282431761293
$ws.Range("D2").Value = -18.54418339531084
$ws.Range("E2").Value = -18.54418339531084
$ws.Range("F2").Value = -18.54418339531084
$ws.Range("G2").Value = -18.54418339531084
$ws.Range("H2").Value = -18.54418339531084
$ws.Range("I2").Value = -18.54418339531084
$ws.Range("J2").Value = -18.54418339531084
$ws.Range("K2").Value = -18.54418339531084

$ws.Range("B3").Value = -18.54418339531084
$ws.Range("C3").Value = -18.54418339531084
$ws.Range("D3").Value = -18.54418339531084
$ws.Range("E3").Value = -18.54418339531084
$ws.Range("F3").Value = -18.54418339531084
$ws.Range("G3").Value = -18.54418339531084
$ws.Range("H3").Value = -18.54418339531084
$ws.Range("I3").Value = 2.769488517428675
$ws.Range("J3").Value = -18.54418339531084
$ws.Range("K3").Value = -18.54418339531084

$ws.Range("B4").Value = -18.54418339531084
$ws.Range("C4").Value = 2.347096730138471
$ws.Range("D4").Value = 2.384757799406882
$ws.Range("E4").Value = -18.54418339531084
$ws.Range("F4").Value = -18.54418339531084
$ws.Range("G4").Value = -18.54418339531084
$ws.Range("H4").Value = 1.574672787248206
$ws.Range("I4").Value = -18.54418339531084
$ws.Range("J4").Value = 2.244841529741674
$ws.Range("K4").Value = -18.54418339531084

$ws.Range("B5").Value = -18.54418339531084
$ws.Range("C5").Value = 1.893459984842948
$ws.Range("D5").Value = -18.54418339531084
$ws.Range("E5").Value = -18.54418339531084
$ws.Range("F5").Value = -18.54418339531084
$ws.Range("G5").Value = 2.904086698956681
$ws.Range("H5").Value = -18.54418339531084
$ws.Range("I5").Value = -18.54418339531084
$ws.Range("J5").Value = -18.54418339531084
$ws.Range("K5").Value = -18.54418339531084

$ws.Range("B6").Value = -18.54418339531084
$ws.Range("C6").Value = -18.54418339531084
$ws.Range("D6").Value = -18.54418339531084
$ws.Range("E6").Value = -18.54418339531084
$ws.Range("F6").Value = -18.54418339531084
$ws.Range("G6").Value = -18.54418339531084
$ws.Range("H6").Value = -18.54418339531084
$ws.Range("I6").Value = -18.54418339531084
$ws.Range("J6").Value = -18.54418339531084
$ws.Range("K6").Value = -18.54418339531084

$ws.Range("B7").Value = 2.622132338783376
$ws.Range("C7").Value = -18.54418339531084
$ws.Range("D7").Value = -18.54418339531084
$ws.Range("E7").Value = -18.54418339531084
$ws.Range("F7").Value = -18.54418339531084
$ws.Range("G7").Value = -18.54418339531084
$ws.Range("H7").Value = -18.54418339531084
$ws.Range("I7").Value = -18.54418339531084
$ws.Range("J7").Value = -18.54418339531084
$ws.Range("K7").Value = -18.54418339531084

$ws.Range("B8").Value = -18.54418339531084
$ws.Range("C8").Value = -18.54418339531084
$ws.Range("D8").Value = -18.54418339531084
$ws.Range("E8").Value = 1.704366162930887
$ws.Range("F8").Value = -18.54418339531084
$ws.Range("G8").Value = -18.54418339531084
$ws.Range("H8").Value = -18.54418339531084
$ws.Range("I8").Value = -18.54418339531084
$ws.Range("J8").Value = -18.54418339531084
$ws.Range("K8").Value = -18.54418339531084

$ws.Range("B9").Value = 3.79112217317202
$ws.Range("C9").Value = -18.54418339531084
$ws.Range("D9").Value = -18.54418339531084
$ws.Range("E9").Value = -18.54418339531084
$ws.Range("F9").Value = -18.54418339531084
$ws.Range("G9").Value = -18.54418339531084
$ws.Range("H9").Value = -18.54418339531084
$ws.Range("I9").Value = -18.54418339531084
$ws.Range("J9").Value = -18.54418339531084
$ws.Range("K9").Value = -18.54418339531084

$ws.Range("B10").Value = -18.54418339531084
$ws.Range("C10").Value = -18.54418339531084
$ws.Range("D10").Value = -18.54418339531084
$ws.Range("E10").Value = -18.54418339531084
$ws.Range("F10").Value = -18.54418339531084
$ws.Range("G10").Value = -18.54418339531084
$ws.Range("H10").Value = -18.54418339531084
$ws.Range("I10").Value = 1.317329142308572
$ws.Range("J10").Value = -18.54418339531084
$ws.Range("K10").Value = 1.875201803033951

$ws.Range("B11").Value = -18.54418339531084
$ws.Range("C11").Value = -18.54418339531084
$ws.Range("D11").Value = -18.54418339531084
$ws.Range("E11").Value = 2.722846134172103
$ws.Range("F11").Value = -18.54418339531084
$ws.Range("G11").Value = 2.68252101090932
$ws.Range("H11").Value = -18.54418339531084
$ws.Range("I11").Value = -18.54418339531084
$ws.Range("J11").Value = -18.54418339531084
$ws.Range("K11").Value = 1.747424632013583

$ws.Range("B12").Value = -18.54418339531084
$ws.Range("C12").Value = -18.54418339531084
$ws.Range("D12").Value = -18.54418339531084
$ws.Range("E12").Value = -18.54418339531084
$ws.Range("F12").Value = -18.54418339531084
$ws.Range("G12").Value = -18.54418339531084
$ws.Range("H12").Value = -18.54418339531084
$ws.Range("I12").Value = -18.54418339531084
$ws.Range("J12").Value = -18.54418339531084
$ws.Range("K12").Value = -18.54418339531084

$ws.Range("B13").Value = -18.54418339531084
$ws.Range("C13").Value = -18.54418339531084
$ws.Range("D13").Value = -18.54418339531084
$ws.Range("E13").Value = 2.421192960381676
$ws.Range("F13").Value = -18.54418339531084
$ws.Range("G13").Value = -18.54418339531084
$ws.Range("H13").Value = -18.54418339531084
$ws.Range("I13").Value = -18.54418339531084
$ws.Range("J13").Value = 2.072115058560448
$ws.Range("K13").Value = 1.916837343044237

$ws.Range("B14").Value = -18.54418339531084
$ws.Range("C14").Value = -18.54418339531084
$ws.Range("D14").Value = 1.460246947587812
$ws.Range("E14").Value = -18.54418339531084
$ws.Range("F14").Value = -18.54418339531084
$ws.Range("G14").Value = -18.54418339531084
$ws.Range("H14").Value = -18.54418339531084
$ws.Range("I14").Value = -18.54418339531084
$ws.Range("J14").Value = -18.54418339531084
$ws.Range("K14").Value = 2.146422518687583

$ws.Range("B15").Value = -18.54418339531084
$ws.Range("C15").Value = -18.54418339531084
$ws.Range("D15").Value = 1.622280070386915
$ws.Range("E15").Value = -18.54418339531084
$ws.Range("F15").Value = -18.54418339531084
$ws.Range("G15").Value = -18.54418339531084
$ws.Range("H15").Value = -18.54418339531084
$ws.Range("I15").Value = -18.54418339531084
$ws.Range("J15").Value = -18.54418339531084
$ws.Range("K15").Value = -18.54418339531084

$ws.Range("B16").Value = -18.54418339531084
$ws.Range("C16").Value = -18.54418339531084
$ws.Range("D16").Value = -18.54418339531084
$ws.Range("E16").Value = -18.54418339531084
$ws.Range("F16").Value = -18.54418339531084
$ws.Range("G16").Value = -18.54418339531084
$ws.Range("H16").Value = -18.54418339531084
$ws.Range("I16").Value = -18.54418339531084
$ws.Range("J16").Value = 2.261577005527877
$ws.Range("K16").Value = -18.54418339531084

$ws.Range("B17").Value = -18.54418339531084
$ws.Range("C17").Value = 1.115960517905431
$ws.Range("D17").Value = 1.375507759099568
$ws.Range("E17").Value = -18.54418339531084
$ws.Range("F17").Value = -18.54418339531084
$ws.Range("G17").Value = -18.54418339531084
$ws.Range("H17").Value = 1.15186891619325
$ws.Range("I17").Value = 1.783940363306195
$ws.Range("J17").Value = 1.839391181398693
$ws.Range("K17").Value = -18.54418339531084

$ws.Range("B18").Value = -18.54418339531084
$ws.Range("C18").Value = -18.54418339531084
$ws.Range("D18").Value = -18.54418339531084
$ws.Range("E18").Value = -18.54418339531084
$ws.Range("F18").Value = -18.54418339531084
$ws.Range("G18").Value = -18.54418339531084
$ws.Range("H18").Value = 1.37670005590021
$ws.Range("I18").Value = 1.121036267084936
$ws.Range("J18").Value = 1.422954189912804
$ws.Range("K18").Value = -18.54418339531084

$ws.Range("B19").Value = -18.54418339531084
$ws.Range("C19").Value = -18.54418339531084
$ws.Range("D19").Value = 1.749244470876364
$ws.Range("E19").Value = -18.54418339531084
$ws.Range("F19").Value = -18.54418339531084
$ws.Range("G19").Value = -18.54418339531084
$ws.Range("H19").Value = 1.661606275211743
$ws.Range("I19").Value = 1.589045379808409
$ws.Range("J19").Value = -18.54418339531084
$ws.Range("K19").Value = -18.54418339531084

$ws.Range("B20").Value = -18.54418339531084
$ws.Range("C20").Value = 1.21149475556595
$ws.Range("D20").Value = 1.580539974259546
$ws.Range("E20").Value = -18.54418339531084
$ws.Range("F20").Value = 4.321924509440943
$ws.Range("G20").Value = -18.54418339531084
$ws.Range("H20").Value = 2.07979540962351
$ws.Range("I20").Value = 1.044093871685079
$ws.Range("J20").Value = -18.54418339531084
$ws.Range("K20").Value = 2.254175279453151

$ws.Range("B21").Value = -18.54418339531084
$ws.Range("C21").Value = 1.421884046818847
$ws.Range("D21").Value = -18.54418339531084
$ws.Range("E21").Value = 2.257983467977059
$ws.Range("F21").Value = -18.54418339531084
$ws.Range("G21").Value = 2.607567136703986
$ws.Range("H21").Value = 2.26644864123069
$ws.Range("I21").Value = -18.54418339531084
$ws.Range("J21").Value = -18.54418339531084
$ws.Range("K21").Value = -18.54418339531084

